$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Santiago / Chile" row (row 2).
$ws.Rows.Item(2).Delete()

# Delete the "New York / Ithaca" row (now row 4, was row 5).
$ws.Rows.Item(4).Delete()

# Match the saved selection state from the edit: whole-row selection on row 4.
$ws.Range("A4:XFD4").Select()
